$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$metadata = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/age-gender-group"

# Version: 7.0.0 -> 8.0.0
$metadata.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$metadata.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet "Include from Age and Gender G" ---
$include = $wb.Worksheets.Item("Include from Age and Gender G")

# System URI: ibm.com -> linuxforhealth.org
$include.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/age-gender-group"

$wb.Save()
